# Update res_line/pl_mw.xlsx data for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.511635313529041
$ws.Range("C2").Value = 0.3057533169054523
$ws.Range("D2").Value = 0.2237948784957311
$ws.Range("F2").Value = 1.160858324775496
$ws.Range("G2").Value = 0.5844827284537999
$ws.Range("H2").Value = 0.7215210258933169
$ws.Range("J2").Value = 0.1979572306752662
$ws.Range("L2").Value = 0.3800416990446536
$ws.Range("O2").Value = 2.591176006116086

# Row 3
$ws.Range("B3").Value = 1.37173675491232
$ws.Range("C3").Value = 0.2928790283801561
$ws.Range("D3").Value = 0.2228563371044814
$ws.Range("F3").Value = 1.171074980692978
$ws.Range("G3").Value = 0.5909563955190933
$ws.Range("H3").Value = 0.7291739916712459
$ws.Range("J3").Value = 0.2005902604390526
$ws.Range("L3").Value = 0.3690196140506288
$ws.Range("O3").Value = 2.620626958114826

# Row 4
$ws.Range("B4").Value = 1.285719270385414
$ws.Range("C4").Value = 0.2849978495569019
$ws.Range("D4").Value = 0.2223529202767267
$ws.Range("F4").Value = 1.178106582284876
$ws.Range("G4").Value = 0.5954368468716922
$ws.Range("H4").Value = 0.7342630972591238
$ws.Range("J4").Value = 0.2023017393702178
$ws.Range("L4").Value = 0.3623503234905598
$ws.Range("O4").Value = 2.64058839764246

# Row 5
$ws.Range("B5").Value = 1.250638978073539
$ws.Range("C5").Value = 0.2817924790285815
$ws.Range("D5").Value = 0.2221661614860651
$ws.Range("F5").Value = 1.181162684761148
$ws.Range("G5").Value = 0.5973895945460299
$ws.Range("H5").Value = 0.7364350552333008
$ws.Range("J5").Value = 0.2030230243065096
$ws.Range("L5").Value = 0.3596574962721917
$ws.Range("O5").Value = 2.64919474363117

# Row 6
$ws.Range("B6").Value = 1.24481233734258
$ws.Range("C6").Value = 0.2812606180315811
$ws.Range("D6").Value = 0.2221362629196832
$ws.Range("F6").Value = 1.181681661995611
$ws.Range("G6").Value = 0.5977215067837633
$ws.Range("H6").Value = 0.7368016325330657
$ws.Range("J6").Value = 0.2031442334424343
$ws.Range("L6").Value = 0.3592118691726967
$ws.Range("O6").Value = 2.650652309002339

# Row 7
$ws.Range("B7").Value = 1.285246273180121
$ws.Range("C7").Value = 0.2849545949821106
$ws.Range("D7").Value = 0.2223503270401181
$ws.Range("F7").Value = 1.17814702604165
$ws.Range("G7").Value = 0.5954626686625559
$ws.Range("H7").Value = 0.7342919918183597
$ws.Range("J7").Value = 0.2023113703191672
$ws.Range("L7").Value = 0.3623139056862783
$ws.Range("O7").Value = 2.640702555742195

# Row 8
$ws.Range("B8").Value = 1.463424485558221
$ws.Range("C8").Value = 0.3013095314346117
$ws.Range("D8").Value = 0.2234561895274254
$ws.Range("F8").Value = 1.164223550652274
$ws.Range("G8").Value = 0.5866097877977694
$ws.Range("H8").Value = 0.7240787917132749
$ws.Range("J8").Value = 0.1988454312525434
$ws.Range("L8").Value = 0.3762210192579261
$ws.Range("O8").Value = 2.60094055489914

# Row 9
$ws.Range("B9").Value = 1.811793079356107
$ws.Range("C9").Value = 0.3335572675707681
$ws.Range("D9").Value = 0.2262000979472418
$ws.Range("F9").Value = 1.142941980095166
$ws.Range("G9").Value = 0.5732705513715004
$ws.Range("H9").Value = 0.7071464843319433
$ws.Range("J9").Value = 0.1928003845791233
$ws.Range("L9").Value = 0.4042639368831686
$ws.Range("O9").Value = 2.53789265988722

# Row 10
$ws.Range("B10").Value = 2.067008022919481
$ws.Range("C10").Value = 0.3573430506120587
$ws.Range("D10").Value = 0.2285635291721348
$ws.Range("F10").Value = 1.130983576691804
$ws.Range("G10").Value = 0.5659351905558481
$ws.Range("H10").Value = 0.6965936454446222
$ws.Range("J10").Value = 0.1888165719355266
$ws.Range("L10").Value = 0.4253274838769698
$ws.Range("O10").Value = 2.500699473837585

# Row 11
$ws.Range("B11").Value = 2.182934074001139
$ws.Range("C11").Value = 0.3681813882657252
$ws.Range("D11").Value = 0.2297135085234743
$ws.Range("F11").Value = 1.126343138008579
$ws.Range("G11").Value = 0.5631362984456274
$ws.Range("H11").Value = 0.6922027297454036
$ws.Range("J11").Value = 0.1871034095587385
$ws.Range("L11").Value = 0.4350079246388532
$ws.Range("O11").Value = 2.485767726131144

# Row 12
$ws.Range("B12").Value = 2.226805318995616
$ws.Range("C12").Value = 0.3722879081167321
$ws.Range("D12").Value = 0.2301596751013335
$ws.Range("F12").Value = 1.124700988383601
$ws.Range("G12").Value = 0.5621540255735198
$ws.Range("H12").Value = 0.6905989163730908
$ws.Range("J12").Value = 0.1864689209466119
$ws.Range("L12").Value = 0.4386876162703288
$ws.Range("O12").Value = 2.480399793877268

# Row 13
$ws.Range("B13").Value = 2.217358123209181
$ws.Range("C13").Value = 0.3714033994131398
$ws.Range("D13").Value = 0.230063110590919
$ws.Range("F13").Value = 1.125049534463081
$ws.Range("G13").Value = 0.5623621201779656
$ws.Range("H13").Value = 0.6909417048627517
$ws.Range("J13").Value = 0.1866049357141852
$ws.Range("L13").Value = 0.4378945137535197
$ws.Range("O13").Value = 2.481543127394247

# Row 14
$ws.Range("B14").Value = 2.186543952539068
$ws.Range("C14").Value = 0.3685191903245766
$ws.Range("D14").Value = 0.2297500009740929
$ws.Range("F14").Value = 1.126205730257006
$ws.Range("G14").Value = 0.5630539295237753
$ws.Range("H14").Value = 0.6920696018288339
$ws.Range("J14").Value = 0.1870509242355824
$ws.Range("L14").Value = 0.4353103775541456
$ws.Range("O14").Value = 2.485320359998667

# Row 15
$ws.Range("B15").Value = 2.167665732830415
$ws.Range("C15").Value = 0.3667528159933795
$ws.Range("D15").Value = 0.2295596030012206
$ws.Range("F15").Value = 1.126928924309674
$ws.Range("G15").Value = 0.5634877966030274
$ws.Range("H15").Value = 0.6927681466770679
$ws.Range("J15").Value = 0.1873259606752704
$ws.Range("L15").Value = 0.4337293247784118
$ws.Range("O15").Value = 2.487671338432875

# Row 16
$ws.Range("B16").Value = 2.05942828639229
$ws.Range("C16").Value = 0.3566350757371026
$ws.Range("D16").Value = 0.2284898755310394
$ws.Range("F16").Value = 1.131302938261776
$ws.Range("G16").Value = 0.5661289514994081
$ws.Range("H16").Value = 0.6968888471618868
$ws.Range("J16").Value = 0.1889305250796998
$ws.Range("L16").Value = 0.4246968055900453
$ws.Range("O16").Value = 2.501715344273407

# Row 17
$ws.Range("B17").Value = 1.992982058609527
$ws.Range("C17").Value = 0.3504325670673722
$ws.Range("D17").Value = 0.2278527529793024
$ws.Range("F17").Value = 1.134191108683616
$ws.Range("G17").Value = 0.5678871828113685
$ws.Range("H17").Value = 0.699521697679117
$ws.Range("J17").Value = 0.1899402531180119
$ws.Range("L17").Value = 0.4191807050815726
$ws.Range("O17").Value = 2.510840377596566

# Row 18
$ws.Range("B18").Value = 1.954747864378874
$ws.Range("C18").Value = 0.3468667675080894
$ws.Range("D18").Value = 0.2274933455312151
$ws.Range("F18").Value = 1.135927553675494
$ws.Range("G18").Value = 0.5689490956876639
$ws.Range("H18").Value = 0.7010745977291606
$ws.Range("J18").Value = 0.190530349213808
$ws.Range("L18").Value = 0.416017278671049
$ws.Range("O18").Value = 2.516275905190881

# Row 19
$ws.Range("B19").Value = 1.941799741610055
$ws.Range("C19").Value = 0.3456597542017903
$ws.Range("D19").Value = 0.2273728690197458
$ws.Range("F19").Value = 1.136528404058971
$ws.Range("G19").Value = 0.5693173285517119
$ws.Range("H19").Value = 0.7016070038553153
$ws.Range("J19").Value = 0.1907317476001849
$ws.Range("L19").Value = 0.4149477998065407
$ws.Range("O19").Value = 2.518148389697316

# Row 20
$ws.Range("B20").Value = 2.00005705646646
$ws.Range("C20").Value = 0.3510926590940642
$ws.Range("D20").Value = 0.2279198467527266
$ws.Range("F20").Value = 1.133875869506056
$ws.Range("G20").Value = 0.5676947748499472
$ws.Range("H20").Value = 0.6992374358427185
$ws.Range("J20").Value = 0.1898318005164858
$ws.Range("L20").Value = 0.4197669438622427
$ws.Range("O20").Value = 2.509849639057393

# Row 21
$ws.Range("B21").Value = 2.195595586818456
$ws.Range("C21").Value = 0.3693662932022335
$ws.Range("D21").Value = 0.2298416791577154
$ws.Range("F21").Value = 1.125863003262339
$ws.Range("G21").Value = 0.562848620194913
$ws.Range("H21").Value = 0.6917367115560324
$ws.Range("J21").Value = 0.1869195398934274
$ws.Range("L21").Value = 0.4360690252246258
$ws.Range("O21").Value = 2.48420311790764

# Row 22
$ws.Range("B22").Value = 2.323230205707205
$ws.Range("C22").Value = 0.381322292909573
$ws.Range("D22").Value = 0.2311600112330581
$ws.Range("F22").Value = 1.12129698137803
$ws.Range("G22").Value = 0.5601338413679571
$ws.Range("H22").Value = 0.6871780629456339
$ws.Range("J22").Value = 0.1850992618383627
$ws.Range("L22").Value = 0.4468043924625249
$ws.Range("O22").Value = 2.469111237569365

# Row 23
$ws.Range("B23").Value = 2.255124812316751
$ws.Range("C23").Value = 0.3749400551896827
$ws.Range("D23").Value = 0.2304507142096952
$ws.Range("F23").Value = 1.123672532274981
$ws.Range("G23").Value = 0.5615412896915899
$ws.Range("H23").Value = 0.6895796593637016
$ws.Range("J23").Value = 0.1860631794098273
$ws.Range("L23").Value = 0.4410673941246586
$ws.Range("O23").Value = 2.477013099057217

# Row 24
$ws.Range("B24").Value = 1.996858553689265
$ws.Range("C24").Value = 0.3507942311369732
$ws.Range("D24").Value = 0.2278894922146293
$ws.Range("F24").Value = 1.134018152420275
$ws.Range("G24").Value = 0.5677816033898537
$ws.Range("H24").Value = 0.6993658283925299
$ws.Range("J24").Value = 0.1898808020653462
$ws.Range("L24").Value = 0.4195018808303104
$ws.Range("O24").Value = 2.51029696194783

# Row 25
$ws.Range("B25").Value = 1.717671570986965
$ws.Range("C25").Value = 0.3248160237861271
$ws.Range("D25").Value = 0.2253965428951048
$ws.Range("F25").Value = 1.148053872798378
$ws.Range("G25").Value = 0.5764472628299941
$ws.Range("H25").Value = 0.7113956926338361
$ws.Range("J25").Value = 0.1943553378215883
$ws.Range("L25").Value = 0.3965960373452333
$ws.Range("O25").Value = 2.553347893685881
